$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet1 "Registro de usuarios" - update a few cell values
# ---------------------------------------------------------------------------
$ws1.Range("E2").Value = "sadf"
$ws1.Range("F2").Value = "adsf"
$ws1.Range("E3").Value = "dsaf"
$ws1.Range("F3").Value = "asdf"

# ---------------------------------------------------------------------------
# Sheet2 "Registro de empresas" - rebuild the header row & add lots of data
# ---------------------------------------------------------------------------

# Header row (row 1)
$ws2.Range("C1").Value = "Url web"
$ws2.Range("I1").Value = "Creada? (no rellenar)"
$ws2.Range("M1").Value = "Tipos de empresas"
$ws2.Range("N1").Value = "valor"

# M1/N1 pick up the same bold+fill header formatting used by A1:I1
$ws2.Range("A1").Copy()
$ws2.Range("M1:N1").PasteSpecial(-4122)

# Row 2
$ws2.Range("C2").Value = "sdaf"
$ws2.Range("E2").Value = "asdf"
$ws2.Range("F2").Value = "sdaf"
$ws2.Range("G2").Value = "sdaf"
$ws2.Range("H2").Value = "dsf"
$ws2.Range("I2").Value = "sadf"
$ws2.Range("M2").Value = "Grupo de Investigación de Universidad"
$ws2.Range("N2").Value = 0

# Row 3
$ws2.Range("C3").Value = "asdf"
$ws2.Range("F3").Value = "sdf"
$ws2.Range("G3").Value = "sdf"
$ws2.Range("H3").Value = "sdaf"
$ws2.Range("I3").Value = "dsf"
$ws2.Range("M3").Value = "Centro de I+D+i"
$ws2.Range("N3").Value = 1

# Row 4
$ws2.Range("I4").Value = "asdf"
$ws2.Range("M4").Value = "Desarrollo de software"
$ws2.Range("N4").Value = 2

# Row 5
$ws2.Range("I5").Value = "sadf"
$ws2.Range("M5").Value = "Fabricante de componentes"
$ws2.Range("N5").Value = 3

# Row 6
$ws2.Range("I6").Value = "sadf"
$ws2.Range("M6").Value = "Fabricante de módulos"
$ws2.Range("N6").Value = 4

# Row 7
$ws2.Range("M7").Value = "Fabricante de sistemas"
$ws2.Range("N7").Value = 5

# Row 8
$ws2.Range("M8").Value = "Ingeniería "
$ws2.Range("N8").Value = 6

# Row 9
$ws2.Range("M9").Value = "Distribución de productos"
$ws2.Range("N9").Value = 7

# Row 10
$ws2.Range("M10").Value = "Consultoría de I+D+i"
$ws2.Range("N10").Value = 8

# Row 11
$ws2.Range("M11").Value = "startup"
$ws2.Range("N11").Value = 9

# Row 12
$ws2.Range("M12").Value = "Aceleradora"
$ws2.Range("N12").Value = 10

# Row 13
$ws2.Range("M13").Value = "Incubadora"
$ws2.Range("N13").Value = 11

# Row 14
$ws2.Range("M14").Value = "Venture capital"
$ws2.Range("N14").Value = 12

# Row 15
$ws2.Range("M15").Value = "Business Angel"
$ws2.Range("N15").Value = 13

# Row 16
$ws2.Range("M16").Value = "Corporate"
$ws2.Range("N16").Value = 14

# Row 17
$ws2.Range("M17").Value = "Empresa industrial usuaria de tecnología"
$ws2.Range("N17").Value = 15

# Row 18
$ws2.Range("M18").Value = "Hospital o centro sanitario"
$ws2.Range("N18").Value = 16

# Row 19
$ws2.Range("M19").Value = "Medio de comunicación"
$ws2.Range("N19").Value = 17

# Row 20
$ws2.Range("M20").Value = "Empresa de servicios"
$ws2.Range("N20").Value = 18

# Row 21
$ws2.Range("M21").Value = "Administración pública"
$ws2.Range("N21").Value = 19

# ---------------------------------------------------------------------------
# Column widths on sheet2 (new columns A,B,C,E,F,G,H,M get custom widths)
# ColumnWidth is expressed in characters; the stored OOXML width ends up
# being (roughly) ColumnWidth + 5/6, rounded to the nearest 1/6th.
# ---------------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 30.166666666666668
$ws2.Columns.Item(2).ColumnWidth = 40.451822916666664
$ws2.Columns.Item(3).ColumnWidth = 36.877604166666664
$ws2.Columns.Item(5).ColumnWidth = 35.592447916666664
$ws2.Columns.Item(6).ColumnWidth = 50.166666666666664
$ws2.Columns.Item(7).ColumnWidth = 38.592447916666664
$ws2.Columns.Item(8).ColumnWidth = 31.166666666666668
$ws2.Columns.Item(13).ColumnWidth = 42.307291666666664

# ---------------------------------------------------------------------------
# Sheet views: sheet1 becomes the active/selected tab, sheet2 loses it and
# scrolls so column C is the left-most visible column.
# ---------------------------------------------------------------------------
$ws2.Range("I21").Select()
$ws1.Activate()
$ws1.Range("F17").Select()

$app = $ws2.Application
$app.ActiveWindow.ScrollColumn = 3
